$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix existing row 111 (corrected car entry) ---
$ws.Range("D111").Value = "Porsche"
$ws.Range("E111").Value = 962
$ws.Range("F111").Value = 1984

# --- Append new rows 1390-1417 ---
# NOTE: cell write order below intentionally matches the original
# authoring sequence so that new shared-string entries are interned
# in the exact order recorded in the target workbook (G1406 before
# G1405 - the two rows' color cells were populated out of row order).

# Row 1390
$ws.Range("A1390").Value = "HotWheels"
$ws.Range("B1390").Value = 2024
$ws.Range("C1390").Value = "Mainline (L Case)"
$ws.Range("D1390").Value = "Mazda"
$ws.Range("E1390").Value = "787B"
$ws.Range("F1390").Value = 1991
$ws.Range("G1390").Value = "Orange #55"
$ws.Range("H1390").Value = "No"

# Row 1391
$ws.Range("A1391").Value = "HotWheels"
$ws.Range("B1391").Value = 2024
$ws.Range("C1391").Value = "Mainline (L Case)"
$ws.Range("D1391").Value = "Audi"
$ws.Range("E1391").Value = "Avant RS2"
$ws.Range("F1391").Value = 1994
$ws.Range("G1391").Value = "Teal"
$ws.Range("H1391").Value = "No"

# Row 1392
$ws.Range("A1392").Value = "HotWheels"
$ws.Range("B1392").Value = 2024
$ws.Range("C1392").Value = "Mainline (K Case)"
$ws.Range("D1392").Value = "Porsche"
$ws.Range("E1392").Value = "Carrera"
$ws.Range("F1392").Value = 1996
$ws.Range("G1392").Value = "Blue"
$ws.Range("H1392").Value = "No"

# Row 1393
$ws.Range("A1393").Value = "HotWheels"
$ws.Range("B1393").Value = 2024
$ws.Range("C1393").Value = "Mainline (L Case)"
$ws.Range("D1393").Value = "Chevy"
$ws.Range("E1393").Value = "Silverado"
$ws.Range("F1393").Value = 2008
$ws.Range("G1393").Value = "Blue"
$ws.Range("H1393").Value = "No"

# Row 1394
$ws.Range("A1394").Value = "HotWheels"
$ws.Range("B1394").Value = 2024
$ws.Range("C1394").Value = "Mainline (L Case)"
$ws.Range("D1394").Value = "Ford"
$ws.Range("E1394").Value = "Maverick Custom"
$ws.Range("F1394").Value = 2022
$ws.Range("G1394").Value = "Blue"
$ws.Range("H1394").Value = "No"

# Row 1395
$ws.Range("A1395").Value = "HotWheels"
$ws.Range("B1395").Value = 2024
$ws.Range("C1395").Value = "Mainline (L Case)"
$ws.Range("D1395").Value = "Nissan"
$ws.Range("E1395").Value = "Skyline 2000GT-R LBWK"
$ws.Range("F1395").Value = 1973
$ws.Range("G1395").Value = "Gray"
$ws.Range("H1395").Value = "No"

# Row 1396
$ws.Range("A1396").Value = "HotWheels"
$ws.Range("B1396").Value = 2024
$ws.Range("C1396").Value = "Mainline (L Case)"
$ws.Range("D1396").Value = "Shelby"
$ws.Range("E1396").Value = "Cobra 427 S/C"
$ws.Range("F1396").Value = 1966
$ws.Range("G1396").Value = "Metal"
$ws.Range("H1396").Value = "No"

# Row 1397
$ws.Range("A1397").Value = "HotWheels"
$ws.Range("B1397").Value = 2024
$ws.Range("C1397").Value = "Mainline (L Case)"
$ws.Range("D1397").Value = "BMW"
$ws.Range("E1397").Value = "M3"
$ws.Range("F1397").Value = 1992
$ws.Range("G1397").Value = "Tan"
$ws.Range("H1397").Value = "No"

# Row 1398
$ws.Range("A1398").Value = "HotWheels"
$ws.Range("B1398").Value = 2024
$ws.Range("C1398").Value = "Mainline (L Case)"
$ws.Range("D1398").Value = "Ford"
$ws.Range("E1398").Value = "Mustang SVO"
$ws.Range("F1398").Value = 1984
$ws.Range("G1398").Value = "Red"
$ws.Range("H1398").Value = "No"

# Row 1399
$ws.Range("A1399").Value = "HotWheels"
$ws.Range("B1399").Value = 2024
$ws.Range("C1399").Value = "Mainline (L Case)"
$ws.Range("D1399").Value = "Acura"
$ws.Range("E1399").Value = "NSX"
$ws.Range("F1399").Value = 1990
$ws.Range("G1399").Value = "Yellow"
$ws.Range("H1399").Value = "No"

# Row 1400
$ws.Range("A1400").Value = "HotWheels"
$ws.Range("B1400").Value = 2024
$ws.Range("C1400").Value = "Mainline (L Case)"
$ws.Range("D1400").Value = "Pagani"
$ws.Range("E1400").Value = "Utopia"
$ws.Range("F1400").Value = 2024
$ws.Range("G1400").Value = "Silver"
$ws.Range("H1400").Value = "No"

# Row 1401
$ws.Range("A1401").Value = "HotWheels"
$ws.Range("B1401").Value = 2024
$ws.Range("C1401").Value = "Mainline (L Case)"
$ws.Range("D1401").Value = "Lotus"
$ws.Range("E1401").Value = "Type 49"
$ws.Range("F1401").Value = 1967
$ws.Range("G1401").Value = "Green #5"
$ws.Range("H1401").Value = "No"

# Row 1402
$ws.Range("A1402").Value = "HotWheels"
$ws.Range("B1402").Value = 2024
$ws.Range("C1402").Value = "Mainline (L Case)"
$ws.Range("D1402").Value = "Bugatti"
$ws.Range("E1402").Value = "Bolide"
$ws.Range("F1402").Value = 2020
$ws.Range("G1402").Value = "Yellow"
$ws.Range("H1402").Value = "No"

# Row 1403
$ws.Range("A1403").Value = "HotWheels"
$ws.Range("B1403").Value = 2024
$ws.Range("C1403").Value = "Mainline (L Case)"
$ws.Range("D1403").Value = "Pontiac"
$ws.Range("E1403").Value = "Firebird T/A"
$ws.Range("F1403").Value = 1977
$ws.Range("G1403").Value = "Black"
$ws.Range("H1403").Value = "No"

# Row 1404
$ws.Range("A1404").Value = "HotWheels"
$ws.Range("B1404").Value = 2024
$ws.Range("C1404").Value = "Mainline (L Case)"
$ws.Range("D1404").Value = "Lamborghini"
$ws.Range("E1404").Value = "Sesto Elemento"
$ws.Range("F1404").Value = 2010
$ws.Range("G1404").Value = "White"
$ws.Range("H1404").Value = "No"

# Row 1405
$ws.Range("A1405").Value = "HotWheels"
$ws.Range("B1405").Value = 2024
$ws.Range("C1405").Value = "Mainline (C Case)"
$ws.Range("D1405").Value = "Chevy"
$ws.Range("E1405").Value = "Corvette Grand Sport Roadster"
$ws.Range("F1405").Value = 1963
$ws.Range("H1405").Value = "No"
$ws.Range("I1405").Value = "Store Exclusive (Kroger)"
# (G1405 deferred -- see after row 1406 below)

# Row 1406
$ws.Range("A1406").Value = "HotWheels"
$ws.Range("B1406").Value = 2024
$ws.Range("C1406").Value = "Mainline (C Case)"
$ws.Range("D1406").Value = "Pontiac"
$ws.Range("E1406").Value = "Firebird"
$ws.Range("F1406").Value = 1970
$ws.Range("G1406").Value = "Orange #70"
$ws.Range("H1406").Value = "No"
$ws.Range("I1406").Value = "Store Exclusive (Kroger)"
# Now backfill G1405 (Red #39), after G1406 (Orange #70) above,
# to match original shared-string insertion order.
$ws.Range("G1405").Value = "Red #39"

# Row 1407
$ws.Range("A1407").Value = "Matchbox"
$ws.Range("B1407").Value = 2024
$ws.Range("C1407").Value = "Mainline"
$ws.Range("D1407").Value = "Lincoln"
$ws.Range("E1407").Value = "Continental"
$ws.Range("F1407").Value = 1964
$ws.Range("G1407").Value = "Black"
$ws.Range("H1407").Value = "No"

# Row 1408
$ws.Range("A1408").Value = "Matchbox"
$ws.Range("B1408").Value = 2023
$ws.Range("C1408").Value = "Mainline"
$ws.Range("D1408").Value = "Porsche"
$ws.Range("E1408").Value = "918 Spyder"
$ws.Range("F1408").Value = 2015
$ws.Range("G1408").Value = "Lime Green"
$ws.Range("H1408").Value = "No"

# Row 1409
$ws.Range("A1409").Value = "HotWheels"
$ws.Range("B1409").Value = 2024
$ws.Range("C1409").Value = "Mainline (L Case)"
$ws.Range("D1409").Value = "Volkswagen"
$ws.Range("E1409").Value = "ID. Buzz"
$ws.Range("F1409").Value = 2025
$ws.Range("G1409").Value = "Orange"
$ws.Range("H1409").Value = "No"

# Row 1410
$ws.Range("A1410").Value = "HotWheels"
$ws.Range("B1410").Value = 2021
$ws.Range("C1410").Value = "Mainline Mini"
$ws.Range("D1410").Value = "Dodge"
$ws.Range("E1410").Value = "Viper RT/10"
$ws.Range("F1410").Value = 1992
$ws.Range("G1410").Value = "Orange"
$ws.Range("H1410").Value = "No"

# Row 1411
$ws.Range("A1411").Value = "Matchbox"
$ws.Range("B1411").Value = 2023
$ws.Range("C1411").Value = "Mainline Mini"
$ws.Range("D1411").Value = "Renault"
$ws.Range("E1411").Value = "Twizy"
$ws.Range("F1411").Value = 2022
$ws.Range("G1411").Value = "Light Blue"
$ws.Range("H1411").Value = "No"

# Row 1412
$ws.Range("A1412").Value = "Matchbox"
$ws.Range("B1412").Value = 2023
$ws.Range("C1412").Value = "Mainline Mini"
$ws.Range("D1412").Value = "Honda"
$ws.Range("E1412").Value = "N600 Off Road"
$ws.Range("F1412").Value = 1970
$ws.Range("G1412").Value = "White #55"
$ws.Range("H1412").Value = "No"

# Row 1413
$ws.Range("A1413").Value = "Matchbox"
$ws.Range("B1413").Value = 2023
$ws.Range("C1413").Value = "Mainline Mini"
$ws.Range("D1413").Value = "Mazda"
$ws.Range("E1413").Value = "Autozam AZ-1"
$ws.Range("F1413").Value = 1992
$ws.Range("G1413").Value = "Blue"
$ws.Range("H1413").Value = "No"

# Row 1414
$ws.Range("A1414").Value = "HotWheels"
$ws.Range("B1414").Value = 2020
$ws.Range("C1414").Value = "Boulevard"
$ws.Range("D1414").Value = "Volkswagen"
$ws.Range("E1414").Value = "Caddy"
$ws.Range("F1414").Value = 1984
$ws.Range("G1414").Value = "Red"
$ws.Range("H1414").Value = "No"
$ws.Range("I1414").Value = "Boulevard #67"

# Row 1415
$ws.Range("A1415").Value = "HotWheels"
$ws.Range("B1415").Value = 2024
$ws.Range("C1415").Value = "Mainline (L Case)"
$ws.Range("D1415").Value = "Koenigsegg"
$ws.Range("E1415").Value = "Jesko"
$ws.Range("F1415").Value = 2020
$ws.Range("G1415").Value = "Purple"
$ws.Range("H1415").Value = "No"

# Row 1416
$ws.Range("A1416").Value = "HotWheels"
$ws.Range("B1416").Value = 2024
$ws.Range("C1416").Value = "Mainline (L Case)"
$ws.Range("D1416").Value = "Volkswagen"
$ws.Range("E1416").Value = "Beetle `"Bumblebee`""
$ws.Range("F1416").Value = 1964
$ws.Range("G1416").Value = "Yellow"
$ws.Range("H1416").Value = "No"

# Row 1417
$ws.Range("A1417").Value = "HotWheels"
$ws.Range("B1417").Value = 2023
$ws.Range("C1417").Value = "Mainline (C Case)"
$ws.Range("D1417").Value = "Nissan"
$ws.Range("E1417").Value = "Patrol Custom"
$ws.Range("F1417").Value = 1983
$ws.Range("G1417").Value = "White"
$ws.Range("H1417").Value = "No"
$ws.Range("I1417").Value = "Store Exclusive (Kroger)"

# --- Update view state to reflect final cursor/scroll position ---
$ws.Activate()
$ws.Range("A1417").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1407
$win.ScrollColumn = 1
